$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume/1h) are stored as plain text in the sheet.
# Excels COM Value setter auto-coerces plain-number-looking strings (e.g. "1.00",
# "577.17") into real numbers, which would change the cell type away from text and
# lose formatting like trailing zeros. Force text by toggling the cell to the "@"
# (Text) number format for the assignment, then restore the default "Normal" style so
# the cell keeps its original (unstyled) appearance, matching the source data.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '66.926.42'
Set-TextValue $ws.Range('E2') '  +0.09%  '
Set-TextValue $ws.Range('D3') '3.075.66'
Set-TextValue $ws.Range('E3') '  -1.20%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.14%  '
Set-TextValue $ws.Range('D5') '577.17'
Set-TextValue $ws.Range('E5') '  -0.16%  '
Set-TextValue $ws.Range('D6') '167.71'
Set-TextValue $ws.Range('E6') '  -2.81%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('E7') '  +0.08%  '
Set-TextValue $ws.Range('D8') '3.071.91'
Set-TextValue $ws.Range('E8') '  -1.21%  '
Set-TextValue $ws.Range('E9') '  -1.81%  '
Set-TextValue $ws.Range('D10') '6.38'
Set-TextValue $ws.Range('E10') '  -1.09%  '
Set-TextValue $ws.Range('E11') '  -1.70%  '
Set-TextValue $ws.Range('D12') '0.470'
Set-TextValue $ws.Range('E12') '  -2.24%  '
Set-TextValue $ws.Range('D13') '0.0000241'
Set-TextValue $ws.Range('E13') '  -2.10%  '
Set-TextValue $ws.Range('D14') '35.96'
Set-TextValue $ws.Range('E14') '  -3.90%  '
Set-TextValue $ws.Range('E15') '  -1.82%  '
Set-TextValue $ws.Range('D16') '3.589.37'
Set-TextValue $ws.Range('E16') '  -1.08%  '
Set-TextValue $ws.Range('D17') '66.852.93'
Set-TextValue $ws.Range('E17') '  +0.06%  '
Set-TextValue $ws.Range('D18') '7.02'
Set-TextValue $ws.Range('E18') '  -1.69%  '
Set-TextValue $ws.Range('D19') '16.93'
Set-TextValue $ws.Range('E19') '  +2.88%  '
Set-TextValue $ws.Range('D20') '3.077.08'
Set-TextValue $ws.Range('E20') '  -1.04%  '
Set-TextValue $ws.Range('D21') '486.65'
Set-TextValue $ws.Range('E21') '  +1.86%  '
Set-TextValue $ws.Range('D22') '0.689'
Set-TextValue $ws.Range('E22') '  -3.54%  '
Set-TextValue $ws.Range('D23') '7.70'
Set-TextValue $ws.Range('E23') '  -3.88%  '
Set-TextValue $ws.Range('D24') '82.75'
Set-TextValue $ws.Range('E24') '  -1.46%  '
Set-TextValue $ws.Range('D25') '12.82'
Set-TextValue $ws.Range('E25') '  -5.39%  '
Set-TextValue $ws.Range('D26') '2.22'
Set-TextValue $ws.Range('E26') '  -3.77%  '
Set-TextValue $ws.Range('D27') '10.30'
Set-TextValue $ws.Range('E27') '  +2.80%  '
Set-TextValue $ws.Range('E28') '  -0.09%  '
Set-TextValue $ws.Range('D29') '7.82'
Set-TextValue $ws.Range('E29') '  -1.45%  '
Set-TextValue $ws.Range('D30') '2.28'
Set-TextValue $ws.Range('E30') '  -6.02%  '
Set-TextValue $ws.Range('E31') '  -1.82%  '
Set-TextValue $ws.Range('D32') '27.65'
Set-TextValue $ws.Range('E32') '  -3.51%  '
Set-TextValue $ws.Range('E33') '  -2.26%  '
Set-TextValue $ws.Range('D34') '0.0₃0909'
Set-TextValue $ws.Range('E34') '  -3.63%  '
Set-TextValue $ws.Range('D35') '1.00'
Set-TextValue $ws.Range('E35') '  +0.06%  '
Set-TextValue $ws.Range('D36') '5.65'
Set-TextValue $ws.Range('E36') '  -3.57%  '
Set-TextValue $ws.Range('D37') '0.951'
Set-TextValue $ws.Range('E37') '  -2.62%  '
Set-TextValue $ws.Range('D38') '46.26'
Set-TextValue $ws.Range('E38') '  -2.09%  '
Set-TextValue $ws.Range('E39') '  +0.92%  '
Set-TextValue $ws.Range('E40') '  -4.86%  '
Set-TextValue $ws.Range('D41') '0.302'
Set-TextValue $ws.Range('E41') '  -2.67%  '
Set-TextValue $ws.Range('E42') '  -4.00%  '
Set-TextValue $ws.Range('D43') '2.758.36'
Set-TextValue $ws.Range('E43') '  -1.73%  '
Set-TextValue $ws.Range('D44') '370.92'
Set-TextValue $ws.Range('E44') '  -2.34%  '
Set-TextValue $ws.Range('D45') '136.15'
Set-TextValue $ws.Range('E45') '  -0.10%  '
Set-TextValue $ws.Range('D46') '0.0345'
Set-TextValue $ws.Range('E46') '  -3.11%  '
Set-TextValue $ws.Range('E47') '  -2.81%  '
Set-TextValue $ws.Range('E48') '  -0.04%  '
Set-TextValue $ws.Range('D49') '24.48'
Set-TextValue $ws.Range('E49') '  -2.15%  '
Set-TextValue $ws.Range('D50') '2.15'
Set-TextValue $ws.Range('E50') '  -2.35%  '
Set-TextValue $ws.Range('E51') '  -1.89%  '
